$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 20:10"

# Update country data rows (country name and/or numeric columns B:H changed)
# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 1756673
$ws.Range("C4").Value = 10870
$ws.Range("D4").Value = 494722
$ws.Range("E4").Value = 1159215
$ws.Range("G4").Value = 629
$ws.Range("H4").Value = 102736

# Row 7: España -> España
$ws.Range("B7").Value = 284986
$ws.Range("C7").Value = 1137
$ws.Range("E7").Value = 60909
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 27119

# Row 11: Alemania -> Alemania
$ws.Range("B11").Value = 182313
$ws.Range("C11").Value = 418
$ws.Range("E11").Value = 10558
$ws.Range("G11").Value = 22
$ws.Range("H11").Value = 8555

# Row 12: India -> India
$ws.Range("B12").Value = 165362
$ws.Range("C12").Value = 7276
$ws.Range("D12").Value = 70788
$ws.Range("E12").Value = 89864
$ws.Range("G12").Value = 176
$ws.Range("H12").Value = 4710

# Row 13: Turquia -> Turquia
$ws.Range("B13").Value = 160979
$ws.Range("C13").Value = 1182
$ws.Range("D13").Value = 124369
$ws.Range("E13").Value = 32149
$ws.Range("G13").Value = 30
$ws.Range("H13").Value = 4461

# Row 16: Canada -> Canada
$ws.Range("B16").Value = 88468
$ws.Range("C16").Value = 949
$ws.Range("E16").Value = 34829

# Row 34: Irlanda -> Irlanda
$ws.Range("B34").Value = 24841
$ws.Range("C34").Value = 38
$ws.Range("E34").Value = 1113
$ws.Range("G34").Value = 8
$ws.Range("H34").Value = 1639

# Row 40: Egipto -> Egipto
$ws.Range("B40").Value = 20793
$ws.Range("C40").Value = 1127
$ws.Range("D40").Value = 5359
$ws.Range("E40").Value = 14589
$ws.Range("G40").Value = 29
$ws.Range("H40").Value = 845

# Row 41: Rumania -> Rumania
$ws.Range("E41").Value = 4927
$ws.Range("G41").Value = 8
$ws.Range("H41").Value = 1235

# Row 99: Eslovaquia -> Sri Lanka
$ws.Range("A99").Value = "Sri Lanka"
$ws.Range("B99").Value = 1524
$ws.Range("C99").Value = 55
$ws.Range("D99").Value = 745
$ws.Range("E99").Value = 769
$ws.Range("H99").Value = 10

# Row 100: Nueva Zelanda -> Eslovaquia
$ws.Range("A100").Value = "Eslovaquia"
$ws.Range("B100").Value = 1520
$ws.Range("C100").Value = 5
$ws.Range("D100").Value = 1332
$ws.Range("E100").Value = 160
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 28

# Row 101: Sri Lanka -> Nueva Zelanda
$ws.Range("A101").Value = "Nueva Zelanda"
$ws.Range("B101").Value = 1504
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 1474
$ws.Range("E101").Value = 8
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 22

# Row 120: Paraguay -> Paraguay
$ws.Range("B120").Value = 900
$ws.Range("C120").Value = 16
$ws.Range("D120").Value = 402
$ws.Range("E120").Value = 487

# Row 125: Principado de Andorra -> Principado de Andorra
$ws.Range("D125").Value = 681
$ws.Range("E125").Value = 31

# Row 153: Liberia -> Yemen
$ws.Range("A153").Value = "Yemen"
$ws.Range("B153").Value = 278
$ws.Range("C153").Value = 22
$ws.Range("D153").Value = 11
$ws.Range("E153").Value = 210
$ws.Range("G153").Value = 4
$ws.Range("H153").Value = 57

# Row 154: Yemen -> Liberia
$ws.Range("A154").Value = "Liberia"
$ws.Range("B154").Value = 269
$ws.Range("C154").Value = 3
$ws.Range("D154").Value = 144
$ws.Range("E154").Value = 98
$ws.Range("H154").Value = 27

# Row 155: Mozambique -> Mozambique
$ws.Range("B155").Value = 233
$ws.Range("C155").Value = 6
$ws.Range("D155").Value = 82
$ws.Range("E155").Value = 149
$ws.Range("G155").Value = 1
$ws.Range("H155").Value = 2

# Row 197: Curazao -> Fiyi
$ws.Range("A197").Value = "Fiyi"
$ws.Range("D197").Value = 15
$ws.Range("H197").Value = 0

# Row 198: Fiyi -> Curazao
$ws.Range("A198").Value = "Curazao"
$ws.Range("D198").Value = 14
$ws.Range("H198").Value = 1

# Row 199: Santa Lucia -> Nueva Caledonia
$ws.Range("A199").Value = "Nueva Caledonia"

# Row 200: Belice -> Santa Lucia
$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("D200").Value = 18
$ws.Range("H200").Value = 0

# Row 201: Nueva Caledonia -> Belice
$ws.Range("A201").Value = "Belice"
$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2

# Row 210: Montserrat -> Seychelles
$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

# Row 211: Seychelles -> Montserrat
$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# Row 213: Papua Nueva Guinea -> Islas Virgenes Britanicas
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

# Row 214: Islas Virgenes Britanicas -> Papua Nueva Guinea
$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0

# Row 215: San Bartolome -> Bonaire, San Eustaquio y Saba
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"

# Row 216: Bonaire, San Eustaquio y Saba -> San Bartolome
$ws.Range("A216").Value = "San Bartolome"

